$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tipo_Viol")

# Clear the total row (row 12), columns B through I which held SUM formulas
$ws.Range("B12:I12").ClearContents()

# Clear the "Total" label in A12 but keep its formatting/style
$ws.Range("A12").ClearContents()

# Update the selection on the sheet to span A12:I12
$ws.Range("A12:I12").Select()
